$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.05445000000000001
$ws.Range("E2").Value = 0.08965000000000001
$ws.Range("G2").Value = 0.1198830409356725
$ws.Range("H2").Value = 0.1198830409356725
$ws.Range("I2").Value = 0.1188271604938271
$ws.Range("J2").Value = 0.1056733024952206
$ws.Range("K2").Value = 116.3
$ws.Range("L2").Value = 0.09446068875893436
$ws.Range("M2").Value = 48.56
$ws.Range("N2").Value = 0.04806493120855192
$ws.Range("O2").Value = 0.4175408426483233
$ws.Range("P2").Value = 48.56
$ws.Range("Q2").Value = 0.04806493120855192
$ws.Range("R2").Value = 0.4175408426483233
$ws.Range("U2").Value = 593.7
$ws.Range("V2").Value = 0.5876472334950016
$ws.Range("W2").Value = 0.1082708019033663
$ws.Range("X2").Value = 0.04410594435819223
$ws.Range("Y2").Value = 0.06416485754517406
$ws.Range("Z2").Value = 2.659926113163523
$ws.Range("AA2").Value = 0.179646926008999
$ws.Range("AB2").Value = 0.04404330926094469
$ws.Range("AC2").Value = 0.1355932252871089
$ws.Range("AD2").Value = 4.218
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 4.218
$ws.Range("AG2").Value = -589.4820000000001
$ws.Range("AH2").Value = 0.004157639391316862
$ws.Range("AI2").Value = 0.003584546169940461
$ws.Range("AJ2").Value = -1.400800345992805
$ws.Range("AK2").Value = -1.011087136246223
$ws.Range("AL2").Value = 0.091
$ws.Range("AM2").Value = 0.091
$ws.Range("AN2").Value = 0.02775
$ws.Range("AO2").Value = 1607.692307692308
$ws.Range("AP2").Value = -3.87817105263158
$ws.Range("AQ2").Value = 1607.692307692308
# Row 3
$ws.Range("D3").Value = 0.0722
$ws.Range("E3").Value = 0.07629999999999999
$ws.Range("G3").Value = 0.1229566453447051
$ws.Range("H3").Value = 0.1229566453447051
$ws.Range("I3").Value = 0.1432125088841507
$ws.Range("J3").Value = 0.1192755609706569
$ws.Range("K3").Value = 58.3
$ws.Range("L3").Value = 0.103589196872779
$ws.Range("M3").Value = 19.6
$ws.Range("N3").Value = 0.04508856682769727
$ws.Range("O3").Value = 0.3361921097770155
$ws.Range("P3").Value = 19.6
$ws.Range("Q3").Value = 0.04508856682769727
$ws.Range("R3").Value = 0.3361921097770155
$ws.Range("U3").Value = 360.2
$ws.Range("V3").Value = 0.8286174373130895
$ws.Range("W3").Value = 0.1649222065063649
$ws.Range("X3").Value = 0.04411646540795676
$ws.Range("Y3").Value = 0.1208057410984082
$ws.Range("Z3").Value = 31.52941176470584
$ws.Range("AA3").Value = 3.760688275310119
$ws.Range("AB3").Value = 0.04404458485279544
$ws.Range("AC3").Value = 3.716643690457324
$ws.Range("AD3").Value = 1.28
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1.28
$ws.Range("AG3").Value = -358.92
$ws.Range("AH3").Value = 0.002935914491490436
$ws.Range("AI3").Value = 0.003007660134404812
$ws.Range("AJ3").Value = -4.736342042755346
$ws.Range("AK3").Value = -5.48975221780361
$ws.Range("AN3").Value = 0.01549636803874092
$ws.Range("AP3").Value = -4.345278450363197
# Row 4
$ws.Range("D4").Value = 0.07150000000000001
$ws.Range("E4").Value = 0.233
$ws.Range("G4").Value = 0.08937544867193108
$ws.Range("H4").Value = 0.08937544867193108
$ws.Range("I4").Value = 0.08111988513998564
$ws.Range("J4").Value = 0.08062971391913977
$ws.Range("K4").Value = 23.3
$ws.Range("L4").Value = 0.08363244795405599
$ws.Range("M4").Value = 6.47
$ws.Range("N4").Value = 0.03835210432720806
$ws.Range("O4").Value = 0.2776824034334764
$ws.Range("P4").Value = 6.47
$ws.Range("Q4").Value = 0.03835210432720806
$ws.Range("R4").Value = 0.2776824034334764
$ws.Range("U4").Value = 77.90000000000001
$ws.Range("V4").Value = 0.4617664493183166
$ws.Range("W4").Value = 0.1349160393746381
$ws.Range("X4").Value = 0.04409542330842771
$ws.Range("Y4").Value = 0.09082061606621041
$ws.Range("Z4").Value = 2.778276391631266
$ws.Range("AA4").Value = 0.2240116306455288
$ws.Range("AB4").Value = 0.04404203366909394
$ws.Range("AC4").Value = 0.1799695969764349
$ws.Range("AD4").Value = 0.369
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.369
$ws.Range("AG4").Value = -77.53100000000001
$ws.Range("AH4").Value = 0.002182540856100172
$ws.Range("AI4").Value = 0.001850839398301642
$ws.Range("AJ4").Value = -0.8504096787285154
$ws.Range("AK4").Value = -0.6382780791806963
$ws.Range("AL4").Value = 0.08599999999999999
$ws.Range("AM4").Value = 0.08599999999999999
$ws.Range("AN4").Value = 0.01506122448979592
$ws.Range("AO4").Value = 262.7906976744187
$ws.Range("AP4").Value = -3.164530612244898
$ws.Range("AQ4").Value = 262.7906976744187
# Row 5
$ws.Range("D5").Value = 0.0211
$ws.Range("E5").Value = -0.09390000000000001
$ws.Range("G5").Value = 0.1619870410367171
$ws.Range("H5").Value = 0.1619870410367171
$ws.Range("I5").Value = 0.152267818574514
$ws.Range("J5").Value = 0.1348001141040792
$ws.Range("K5").Value = 23.5
$ws.Range("L5").Value = 0.126889848812095
$ws.Range("M5").Value = 12.5
$ws.Range("N5").Value = 0.04729474082482028
$ws.Range("O5").Value = 0.5319148936170213
$ws.Range("P5").Value = 12.5
$ws.Range("Q5").Value = 0.04729474082482028
$ws.Range("R5").Value = 0.5319148936170213
$ws.Range("U5").Value = 97.7
$ws.Range("V5").Value = 0.3696556942867953
$ws.Range("W5").Value = 0.08162556443209448
$ws.Range("X5").Value = 0.04428907441816117
$ws.Range("Y5").Value = 0.03733649001393331
$ws.Range("Z5").Value = 1.003576460387992
$ws.Range("AA5").Value = 0.1352822213724692
$ws.Range("AB5").Value = 0.04406536777468631
$ws.Range("AC5").Value = 0.09121685359778292
$ws.Range("AD5").Value = 2.42
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 2.42
$ws.Range("AG5").Value = -95.28
$ws.Range("AH5").Value = 0.009073185362927413
$ws.Range("AI5").Value = 0.007503410641200544
$ws.Range("AJ5").Value = -0.5637202697905573
$ws.Range("AK5").Value = -0.4238057112356551
$ws.Range("AN5").Value = 0.08259385665529009
$ws.Range("AP5").Value = -3.251877133105802
# Row 6
$ws.Range("D6").Value = 0.0374
$ws.Range("E6").Value = 0.103
$ws.Range("G6").Value = 0.1148582600195504
$ws.Range("H6").Value = 0.1148582600195504
$ws.Range("I6").Value = 0.07282502443792767
$ws.Range("J6").Value = 0.0615453589986697
$ws.Range("K6").Value = 11.2
$ws.Range("L6").Value = 0.05474095796676442
$ws.Range("M6").Value = 9.99
$ws.Range("N6").Value = 0.07005610098176719
$ws.Range("O6").Value = 0.8919642857142858
$ws.Range("P6").Value = 9.99
$ws.Range("Q6").Value = 0.07005610098176719
$ws.Range("R6").Value = 0.8919642857142858
$ws.Range("U6").Value = 57.9
$ws.Range("V6").Value = 0.406030855539972
$ws.Range("W6").Value = 0.05348615090735434
$ws.Range("X6").Value = 0.0440636776795776
$ws.Range("Y6").Value = 0.009422473227776743
$ws.Range("Z6").Value = 1.277137613762625
$ws.Range("AA6").Value = 0.0786018929297251
$ws.Range("AB6").Value = 0.04403817746533634
$ws.Range("AC6").Value = 0.03456371546438876
$ws.Range("AD6").Value = 0.149
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0.149
$ws.Range("AG6").Value = -57.751
$ws.Range("AH6").Value = 0.001043790149142901
$ws.Range("AI6").Value = 0.0006499483094800849
$ws.Range("AJ6").Value = -0.6806326533017478
$ws.Range("AK6").Value = -0.3370372748017205
$ws.Range("AL6").Value = 0.005
$ws.Range("AM6").Value = 0.005
$ws.Range("AN6").Value = 0.009551282051282051
$ws.Range("AO6").Value = 2980
$ws.Range("AP6").Value = -3.701987179487179
$ws.Range("AQ6").Value = 2980
